$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.182.18'
$ws.Range('E2').Value = '  +0.26%  '
$ws.Range('D3').Value = '1.765.86'
$ws.Range('E3').Value = '  +2.87%  '
$ws.Range('D4').Value = '''0.9970'
$ws.Range('E4').Value = '  -0.48%  '
$ws.Range('D5').Value = '''312.42'
$ws.Range('E5').Value = '  +1.66%  '
$ws.Range('D6').Value = '''0.9953'
$ws.Range('E6').Value = '  -0.55%  '
$ws.Range('D7').Value = '''0.5155'
$ws.Range('E7').Value = '  +9.72%  '
$ws.Range('D8').Value = '''0.3604'
$ws.Range('E8').Value = '  +5.24%  '
$ws.Range('D9').Value = '''42.27'
$ws.Range('D10').Value = '''0.07292'
$ws.Range('E10').Value = '  +0.35%  '
$ws.Range('D11').Value = '''1.074'
$ws.Range('E11').Value = '  +3.02%  '
$ws.Range('D12').Value = '''0.9936'
$ws.Range('E12').Value = '  -0.77%  '
$ws.Range('E13').Value = '  +2.79%  '
$ws.Range('D14').Value = '''6.018'
$ws.Range('E14').Value = '  +2.60%  '
$ws.Range('D15').Value = '1.754.01'
$ws.Range('E15').Value = '  +2.05%  '
$ws.Range('D16').Value = '''6.899'
$ws.Range('E16').Value = '  +0.15%  '
$ws.Range('D17').Value = '''87.87'
$ws.Range('E17').Value = '  -1.35%  '
$ws.Range('D18').Value = '''0.00001039'
$ws.Range('E18').Value = '  +0.11%  '
$ws.Range('D19').Value = '''0.06408'
$ws.Range('E19').Value = '  +0.91%  '
$ws.Range('D20').Value = '''0.9952'
$ws.Range('E20').Value = '  -0.51%  '
$ws.Range('D21').Value = '''16.70'
$ws.Range('E21').Value = '  +1.12%  '
$ws.Range('D22').Value = '''5.798'
$ws.Range('E22').Value = '  +3.07%  '
$ws.Range('D23').Value = '27.281.86'
$ws.Range('E23').Value = '  +0.44%  '
$ws.Range('D24').Value = '''11.31'
$ws.Range('E24').Value = '  +4.16%  '
$ws.Range('D25').Value = '''2.052'
$ws.Range('E25').Value = '  -3.02%  '
$ws.Range('D26').Value = '''154.29'
$ws.Range('E26').Value = '  -1.71%  '
$ws.Range('D27').Value = '''20.11'
$ws.Range('E27').Value = '  +3.18%  '
$ws.Range('D28').Value = '1.952.05'
$ws.Range('E28').Value = '  +2.16%  '
$ws.Range('D29').Value = '''2.278'
$ws.Range('E29').Value = '  +8.31%  '
$ws.Range('D30').Value = '''120.66'
$ws.Range('E30').Value = '  +0.86%  '
$ws.Range('D31').Value = '''1.057'
$ws.Range('E31').Value = '  +3.79%  '
$ws.Range('D32').Value = '''0.09608'
$ws.Range('E32').Value = '  +4.88%  '
$ws.Range('D33').Value = '''3.592'
$ws.Range('E33').Value = '  +0.08%  '
$ws.Range('D34').Value = '''5.458'
$ws.Range('E34').Value = '  +2.58%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').Value = '''0.05952'
$ws.Range('E35').Value = '  +2.27%  '
$ws.Range('B36').Value = 'VeChain'
$ws.Range('C36').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D36').Value = '''0.02205'
$ws.Range('E36').Value = '  +0.09%  '
$ws.Range('D37').Value = '''11.13'
$ws.Range('E37').Value = '  +1.44%  '
$ws.Range('B38').Value = 'Algorand'
$ws.Range('C38').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D38').Value = '''0.2018'
$ws.Range('E38').Value = '  +1.07%  '
$ws.Range('B39').Value = 'InternetComputer(DFINITY)'
$ws.Range('C39').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D39').Value = '''4.813'
$ws.Range('E39').Value = '  +1.45%  '
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').Value = '''0.6092'
$ws.Range('E40').Value = '  +3.31%  '
$ws.Range('D41').Value = '''1.426'
$ws.Range('E41').Value = '  +2.40%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').Value = '''7.760'
$ws.Range('E42').Value = '  +4.24%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').Value = '''1.115'
$ws.Range('E43').Value = '  -0.52%  '
$ws.Range('D44').Value = '''13.04'
$ws.Range('E44').Value = '  +3.48%  '
$ws.Range('D45').Value = '''3.610'
$ws.Range('E45').Value = '  +1.62%  '
$ws.Range('D46').Value = '''0.5713'
$ws.Range('E46').Value = '  +1.17%  '
$ws.Range('D47').Value = '''120.92'
$ws.Range('E47').Value = '  +2.61%  '
$ws.Range('D48').Value = '''1.867'
$ws.Range('E48').Value = '  +1.52%  '
$ws.Range('B49').Value = 'EOS'
$ws.Range('C49').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D49').Value = '''1.105'
$ws.Range('E49').Value = '  +1.83%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = '''0.06687'
$ws.Range('E50').Value = '  +0.57%  '
$ws.Range('D51').Value = '''70.14'
$ws.Range('E51').Value = '  +0.26%  '
